$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price/Volume columns for the rows being updated
# so numeric-looking strings (e.g. "286.60", "4.12%") are stored as literal text,
# matching the original inlineStr text cells instead of being parsed as numbers.
$ws.Range("D2:E46").NumberFormat = "@"

$ws.Range("D2").Value = "286.60"
$ws.Range("E2").Value = "4.12%"
$ws.Range("D3").Value = "28.35"
$ws.Range("E3").Value = "4.34%"
$ws.Range("D4").Value = "4.936"
$ws.Range("E4").Value = "1.93%"
$ws.Range("D5").Value = "0.06554"
$ws.Range("E5").Value = "2.44%"
$ws.Range("D6").Value = "7.250"
$ws.Range("E6").Value = "4.66%"
$ws.Range("D7").Value = "1.337"
$ws.Range("E7").Value = "11.06%"
$ws.Range("D8").Value = "0.9175"
$ws.Range("E8").Value = "4.26%"
$ws.Range("D9").Value = "0.1568"
$ws.Range("E9").Value = "3.44%"
$ws.Range("D10").Value = "0.06496"
$ws.Range("E10").Value = "29.69%"
$ws.Range("D11").Value = "0.07696"
$ws.Range("E11").Value = "2.05%"
$ws.Range("D12").Value = "0.02976"
$ws.Range("E12").Value = "0.16%"
$ws.Range("E13").Value = "-0.34%"
$ws.Range("D14").Value = "0.001595"
$ws.Range("E14").Value = "1.42%"
$ws.Range("D15").Value = "0.0006537"
$ws.Range("E15").Value = "1.92%"
$ws.Range("D16").Value = "0.006014"
$ws.Range("E16").Value = "-2.97%"
$ws.Range("D17").Value = "3.487"
$ws.Range("E17").Value = "0.61%"
$ws.Range("D18").Value = "3.378"
$ws.Range("E18").Value = "2.11%"
$ws.Range("D19").Value = "2.240"
$ws.Range("E19").Value = "-1.93%"
$ws.Range("E20").Value = "0.66%"
$ws.Range("D21").Value = "0.1349"
$ws.Range("E21").Value = "-0.81%"
$ws.Range("D22").Value = "4.030"
$ws.Range("E22").Value = "2.50%"
$ws.Range("B23").Value = "ZBToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D23").Value = "0.1549"
$ws.Range("E23").Value = "12.27%"
$ws.Range("B24").Value = "CoinExToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D24").Value = "0.04470"
$ws.Range("E24").Value = "1.10%"
$ws.Range("D25").Value = "0.001186"
$ws.Range("E25").Value = "0.80%"
$ws.Range("D26").Value = "0.004346"
$ws.Range("E26").Value = "12.83%"
$ws.Range("D28").Value = "0.0001180"
$ws.Range("E28").Value = "-1.93%"
$ws.Range("D29").Value = "0.0001635"
$ws.Range("E29").Value = "-15.78%"
$ws.Range("D40").Value = "0.04152"
$ws.Range("E40").Value = "0.54%"
$ws.Range("D41").Value = "0.007048"
$ws.Range("E41").Value = "3.21%"
$ws.Range("E42").Value = "20.32%"
$ws.Range("E43").Value = "-5.34%"
$ws.Range("D44").Value = "0.01249"
$ws.Range("E44").Value = "8.78%"
$ws.Range("D45").Value = "0.00005556"
$ws.Range("E45").Value = "7.62%"
$ws.Range("D46").Value = "1.572"
$ws.Range("E46").Value = "-4.72%"
